# Update the "updated_at" (column F) timestamps for four vendor-selection
# rows to reflect a new incident:
#   - row 4  (MEH_PRD_NAS)                -> 2025-12-25 23:25:54
#   - row 11 (AHM_PRD_NAS)                -> 2025-12-25 23:25:54
#   - row 12 (HOST-3.123.68.65)           -> 2025-12-25 23:25:57
#   - row 16 (PV_RPT_DB.ahm.lambdacro.com)-> 2025-12-25 23:25:57
#
# The order below is chosen so the workbook's shared-string table ends up
# laid out the same way Excel itself produced it: row 12 is the sole owner
# of its old shared string, so touching it first/last lets the engine
# reuse that slot in place instead of always appending new entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F12").Value = "2025-12-25 23:25:54"
$ws.Range("F4").Value  = "2025-12-25 23:25:54"
$ws.Range("F11").Value = "2025-12-25 23:25:54"

$ws.Range("F12").Value = "2025-12-25 23:25:57"
$ws.Range("F16").Value = "2025-12-25 23:25:57"
